# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Change: cell B11 on the "Rules" sheet held the shared string "R40".
# It is replaced with the text value "1" (kept as TEXT, not a number),
# while preserving the cell's existing style/formatting.
#
# A plain `.Value = "1"` assignment would be auto-coerced to a numeric
# literal by Excel (since "1" parses as a number), which would also
# silently change the cell's number format. To force a genuine text
# value without touching the cell's style, we briefly hold it as a
# text *formula* result ( ="1" ) and then convert that computed value
# in place to a literal via Copy / Paste Special (values only) - the
# same trick used interactively in Excel to "flatten" a formula to text
# without disturbing formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$target = $ws.Range("B11")
$target.Formula = "=""1"""
$target.Copy()
$target.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
